$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 377.8421
$ws.Range("I15").Value = 377.8421
$ws.Range("K15").Value = 1133.5263
$ws.Range("M15").Value = -964.5263
$ws.Range("H33").Value = 255.1875
$ws.Range("I33").Value = 270.2143
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 270.2143
$ws.Range("L33").Value = 150
$ws.Range("M33").Value = -41.21429999999998
$ws.Range("N33").Value = -608
$ws.Range("H38").Value = 2402.25
$ws.Range("I38").Value = 63.6
$ws.Range("J38").Value = 6300
$ws.Range("K38").Value = 190.8
$ws.Range("L38").Value = 18900
$ws.Range("M38").Value = 181.2
$ws.Range("N38").Value = -19644
$ws.Range("H43").Value = 9999.666999999999
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").Value = $null
$ws.Range("H88").Value = 6075.75
$ws.Range("I88").Value = 3400
$ws.Range("K88").Value = 3400
$ws.Range("M88").Value = -2994
$ws.Range("H91").Value = 6075.75
$ws.Range("I91").Value = 3400
$ws.Range("K91").Value = 3400
$ws.Range("M91").Value = -1996
$ws.Range("H135").Value = 1789.7
$ws.Range("I135").Value = 1613.25
$ws.Range("K135").Value = 14519.25
$ws.Range("M135").Value = -11984.25
$ws.Range("H137").Value = 6572.1177
$ws.Range("I137").Value = 6148.8335
$ws.Range("J137").Value = 7588
$ws.Range("K137").Value = 18446.5005
$ws.Range("L137").Value = 22764
$ws.Range("M137").Value = -15896.5005
$ws.Range("N137").Value = -27864
$ws.Range("H138").Value = 1488.3334
$ws.Range("I138").Value = 732.5
$ws.Range("K138").Value = 2197.5
$ws.Range("M138").Value = 2942.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1402.75
$ws.Range("I2").Value = 1402.75
$ws.Range("K2").Value = 1402.75
$ws.Range("M2").Value = -1289.75
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = $null
$ws.Range("H32").Value = 314.53845
$ws.Range("I32").Value = 314.53845
$ws.Range("K32").Value = 314.53845
$ws.Range("M32").Value = -27.53845000000001
$ws.Range("H54").Value = 13072
$ws.Range("J54").Value = 13072
$ws.Range("L54").Value = 13072
$ws.Range("N54").Value = -14610
$ws.Range("H116").Value = 1402.75
$ws.Range("I116").Value = 1402.75
$ws.Range("K116").Value = 1402.75
$ws.Range("M116").Value = 891.25
$ws.Range("H132").Value = 7558.5454
$ws.Range("I132").Value = 5455.5
$ws.Range("K132").Value = 16366.5
$ws.Range("M132").Value = -13836.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1402.75
$ws.Range("I3").Value = 1402.75
$ws.Range("K3").Value = 1402.75
$ws.Range("M3").Value = -1288.75
$ws.Range("H5").Value = 241.33333
$ws.Range("I5").Value = 284.6
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 284.6
$ws.Range("L5").Value = 25
$ws.Range("M5").Value = -171.6
$ws.Range("N5").Value = -251
$ws.Range("H22").Value = 473
$ws.Range("I22").Value = 463
$ws.Range("K22").Value = 463
$ws.Range("M22").Value = -290
$ws.Range("H82").Value = 4333.3335
$ws.Range("I82").Value = 4333.3335
$ws.Range("K82").Value = 4333.3335
$ws.Range("M82").Value = -3950.3335
$ws.Range("H85").Value = 4333.3335
$ws.Range("I85").Value = 4333.3335
$ws.Range("K85").Value = 4333.3335
$ws.Range("M85").Value = -3007.3335
$ws.Range("H107").Value = 774.5
$ws.Range("I107").Value = 799.6667
$ws.Range("K107").Value = 799.6667
$ws.Range("M107").Value = 1120.3333
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5442.25
$ws.Range("I31").Value = 3841.9285
$ws.Range("K31").Value = 3841.9285
$ws.Range("M31").Value = -3546.9285
$ws.Range("H34").Value = 5442.25
$ws.Range("I34").Value = 3841.9285
$ws.Range("K34").Value = 3841.9285
$ws.Range("M34").Value = -3639.9285
$ws.Range("H59").Value = 17994.5
$ws.Range("I59").Value = 17994.5
$ws.Range("K59").Value = 17994.5
$ws.Range("M59").Value = -16849.5
$ws.Range("H62").Value = 6025
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("H65").Value = 6025
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("H96").Value = 11641
$ws.Range("J96").Value = 11641
$ws.Range("L96").Value = 11641
$ws.Range("N96").Value = -17133
$ws.Range("H99").Value = 5000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = $null
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060
$ws.Range("H134").Value = 7878.375
$ws.Range("I134").Value = 4402.6
$ws.Range("K134").Value = 13207.8
$ws.Range("M134").Value = -10672.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 513
$ws.Range("I81").Value = 513
$ws.Range("K81").Value = 1539
$ws.Range("M81").Value = -416
$ws.Range("H84").Value = 513
$ws.Range("I84").Value = 513
$ws.Range("K84").Value = 4617
$ws.Range("M84").Value = 999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3297
$ws.Range("J102").Value = 4000
$ws.Range("L102").Value = 4000
$ws.Range("N102").Value = -7244
$ws.Range("H132").Value = 6631.5
$ws.Range("I132").Value = 2920.5
$ws.Range("K132").Value = 8761.5
$ws.Range("M132").Value = -6231.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5339.087
$ws.Range("J46").Value = 5438.048
$ws.Range("L46").Value = 5438.048
$ws.Range("N46").Value = -5814.048
$ws.Range("H75").Value = 17300
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 17300
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 17300
$ws.Range("M75").Value = $null
$ws.Range("N75").Value = -19172
$ws.Range("H78").Value = 17300
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 17300
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 51900
$ws.Range("M78").Value = $null
$ws.Range("N78").Value = -61260
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = $null
$ws.Range("H140").Value = 78331.336
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7945.5713
$ws.Range("I132").Value = 5723.9
$ws.Range("K132").Value = 17171.7
$ws.Range("M132").Value = -14641.7
